$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '37.250.90'
$ws.Range('E2').Value = '  +1.61%  '
$ws.Range('D3').Value = '2.013.80'
$ws.Range('E3').Value = '  +2.57%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue 'D5' '246.97'
Set-TextValue 'D6' '0.627'
$ws.Range('E6').Value = '  +1.74%  '
Set-TextValue 'D7' '59.92'
$ws.Range('E7').Value = '  -1.48%  '
Set-TextValue 'D9' '0.387'
$ws.Range('E9').Value = '  +2.91%  '
$ws.Range('E10').Value = '  +1.47%  '
$ws.Range('E11').Value = '  +0.79%  '
Set-TextValue 'D12' '15.19'
$ws.Range('E12').Value = '  +5.91%  '
Set-TextValue 'D13' '22.41'
$ws.Range('E13').Value = '  +2.23%  '
$ws.Range('D14').Value = '2.305.42'
$ws.Range('E14').Value = '  +2.39%  '
Set-TextValue 'D15' '0.849'
$ws.Range('E15').Value = '  +1.11%  '
Set-TextValue 'D16' '5.48'
$ws.Range('E16').Value = '  +3.23%  '
$ws.Range('D17').Value = '2.017.40'
$ws.Range('E17').Value = '  +2.97%  '
$ws.Range('D18').Value = '37.129.28'
$ws.Range('E18').Value = '  +1.44%  '
Set-TextValue 'D19' '70.38'
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('E21').Value = '  +2.30%  '
Set-TextValue 'D22' '230.98'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('E24').Value = '  +0.53%  '
Set-TextValue 'D25' '2.36'
$ws.Range('E25').Value = '  +0.32%  '
$ws.Range('E26').Value = '  +2.31%  '
Set-TextValue 'D27' '164.59'
Set-TextValue 'D28' '0.139'
$ws.Range('E28').Value = '  -2.68%  '
Set-TextValue 'D29' '19.75'
$ws.Range('E29').Value = '  +1.56%  '
$ws.Range('E30').Value = '  +12.59%  '
$ws.Range('E31').Value = '  +1.39%  '
$ws.Range('E32').Value = '  +0.96%  '
Set-TextValue 'D33' '0.0658'
$ws.Range('E33').Value = '  +6.37%  '
Set-TextValue 'D34' '4.48'
$ws.Range('E34').Value = '  +0.23%  '
Set-TextValue 'D35' '2.46'
$ws.Range('E35').Value = '  +7.83%  '
$ws.Range('E36').Value = '  -2.46%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  +2.04%  '
Set-TextValue 'D39' '5.35'
$ws.Range('E39').Value = '  -4.50%  '
Set-TextValue 'D40' '0.0984'
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('E41').Value = '  +0.83%  '
$ws.Range('E42').Value = '  +0.67%  '
$ws.Range('E43').Value = '  +1.19%  '
Set-TextValue 'D44' '16.71'
$ws.Range('E44').Value = '  +2.69%  '
Set-TextValue 'D45' '92.11'
$ws.Range('E45').Value = '  +3.68%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D46' '1.06'
$ws.Range('E46').Value = '  +2.00%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.377.13'
$ws.Range('E47').Value = '  +0.60%  '
Set-TextValue 'D48' '7.48'
$ws.Range('E48').Value = '  +4.30%  '
Set-TextValue 'D49' '2.09'
$ws.Range('E49').Value = '  +13.09%  '
$ws.Range('E50').Value = '  +0.14%  '
Set-TextValue 'D51' '46.68'
$ws.Range('E51').Value = '  +5.24%  '
